# Weekly Fruit/Vegetable update: insert 3 new rows of data (a new reporting
# week, 2022-03-08) at the top of the "Sandia" (watermelon) data block,
# pushing the existing rows down by three positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("456:458").Insert()

$newRows = @(
    @(456, 44628, "Extra",   250, 2800, 3000, 2900, 2900, "Región Metropolitana"),
    @(457, 44628, "Primera", 430, 2400, 2700, 2550, 2550, "Región Metropolitana"),
    @(458, 44628, "Segunda", 170, 2300, 2300, 2300, 2300, "Región Metropolitana")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value  = 9
    $ws.Cells.Item($r, 2).Value  = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($r, 3).Value  = "Metropolitana"
    $ws.Cells.Item($r, 4).Value  = $row[1]
    $ws.Cells.Item($r, 5).Value  = 13
    $ws.Cells.Item($r, 6).Value  = 100112028
    $ws.Cells.Item($r, 7).Value  = "Sandia"
    $ws.Cells.Item($r, 8).Value  = "Sin especificar"
    $ws.Cells.Item($r, 9).Value  = $row[2]
    $ws.Cells.Item($r, 10).Value = $row[3]
    $ws.Cells.Item($r, 11).Value = $row[4]
    $ws.Cells.Item($r, 12).Value = $row[5]
    $ws.Cells.Item($r, 13).Value = $row[6]
    $ws.Cells.Item($r, 14).Value = "`$/unidad"
    $ws.Cells.Item($r, 15).Value = $row[8]
    $ws.Cells.Item($r, 16).Value = $row[7]
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
